$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Heading: "Tampilan Score " -> "Tampilan Score Akhir Game "
#    (still bold, same run formatting as before)
# ---------------------------------------------------------------------
$headingRange = $d.Content
$foundHeading = $headingRange.Find.Execute(
    "Tampilan Score ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Tampilan Score Akhir Game ", 2)
Write-Host "Heading replace found:" $foundHeading

# ---------------------------------------------------------------------
# 2) Body paragraph describing the "Score Akhir" scene: reword "jenis"
#    to "kategori" and append extra explanatory sentence, with the two
#    new "game" mentions in italics.
# ---------------------------------------------------------------------
$oldBody = "Setelah pemain menjawab semua pertanyaan pada salah satu jenis game yang disediakan maka pemain akan menuju scene Score Akhir di scene ini pemain akan melihat hasil score dari game yang pemain mainkan."
$newBody = "Setelah pemain menjawab semua pertanyaan pada salah satu kategori game yang disediakan maka pemain akan menuju scene Score Akhir di scene ini pemain akan melihat hasil score dari game yang pemain mainkan dari semua kategori game dan jumlah keseluruhan score game."

$bodyRange = $d.Content
$foundBody = $bodyRange.Find.Execute(
    $oldBody, $true, $false, $false, $false, $false,
    $true, 1, $false, $newBody, 2)
Write-Host "Body replace found:" $foundBody

# Italicize the first new "game" occurrence (" dari semua kategori game dan ...")
$searchRange1 = $d.Content
$found1 = $searchRange1.Find.Execute(
    "kategori game dan", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
Write-Host "Locate game#1 found:" $found1
if ($found1) {
    $gameStart1 = $searchRange1.Start + 9
    $gameEnd1 = $gameStart1 + 4
    $gameRange1 = $d.Range($gameStart1, $gameEnd1)
    Write-Host "game#1 text: [" $gameRange1.Text "]"
    $gameRange1.Font.Italic = $true
}

# Italicize the second new "game" occurrence ("...keseluruhan score game.")
$searchRange2 = $d.Content
$found2 = $searchRange2.Find.Execute(
    "score game.", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
Write-Host "Locate game#2 found:" $found2
if ($found2) {
    $gameStart2 = $searchRange2.Start + 6
    $gameEnd2 = $gameStart2 + 4
    $gameRange2 = $d.Range($gameStart2, $gameEnd2)
    Write-Host "game#2 text: [" $gameRange2.Text "]"
    $gameRange2.Font.Italic = $true
}

# ---------------------------------------------------------------------
# 3) Header page-number field cached text: "54" -> "57"
# ---------------------------------------------------------------------
$section = $d.Sections.Item(1)
$primaryHeader = $section.Headers.Item(1)
if ($primaryHeader.Exists -and $primaryHeader.Range.Fields.Count -ge 1) {
    $pageField = $primaryHeader.Range.Fields.Item(1)
    Write-Host "Header field before: [" $pageField.Result.Text "]"
    $pageField.Result.Text = "57"
    Write-Host "Header field after: [" $pageField.Result.Text "]"
}

Write-Host "Done"
